$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)

# Paragraph 1
$tr = $shp.TextFrame.TextRange
$tr.Text = "A rendszer legyen képes az autók mellett a szerelők nyilvántartására is."
$tr.LanguageID = "hu-HU"

# Paragraph 2
$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("`rA rendszer legyen képes részletes adattárolásra (egy tulajdonoshoz több autót is hozzá lehessen rendelni.)")
$ins.LanguageID = "hu-HU"

# Paragraph 3 (three separate runs)
$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("`rLegyünk ")
$ins.LanguageID = "hu-HU"

$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("képesek külön-külön ")
$ins.LanguageID = "hu-HU"

$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("szerkeszteni az ilyen jellegű autókat (melyek egy tulajdonoshoz vannak rendelve.)")
$ins.LanguageID = "hu-HU"

# Paragraph 4
$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("`rA programunk továbbfejlesztett verziója is, egyszerű, letisztult, könnyedén kezelhető, felhasználóbarát felülettel rendelkezzen.")
$ins.LanguageID = "hu-HU"

# Paragraph 5: empty trailing paragraph with no bullet
$cur = $shp.TextFrame.TextRange
$ins = $cur.InsertAfter("`rTEMP")
$ins.LanguageID = "hu-HU"

$all = $shp.TextFrame.TextRange
$count = $all.Paragraphs().Count
$lastPara = $all.Paragraphs($count, 1)
$lastPara.ParagraphFormat.Bullet.Type = 0
$lastPara.Text = ""
